$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows (170, 171) at the end that preserve the ORIGINAL
# (pre-edit) data currently held in rows 168 and 169, before those rows
# get overwritten with their new weekly values below.
$ws.Range("A168:R168").Copy($ws.Range("A170:R170"))
$ws.Range("A169:R169").Copy($ws.Range("A171:R171"))

# Update row 168 with the new weekly values.
$ws.Cells.Item(168, 4).Value2 = 44595
$ws.Cells.Item(168, 10).Value2 = 300
$ws.Cells.Item(168, 11).Value2 = 9000
$ws.Cells.Item(168, 12).Value2 = 10000
$ws.Cells.Item(168, 13).Value2 = 9567
$ws.Cells.Item(168, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(168, 16).Value2 = 191
$ws.Cells.Item(168, 17).Value2 = 50

# Update row 169 with the new weekly values.
$ws.Cells.Item(169, 4).Value2 = 44595
$ws.Cells.Item(169, 10).Value2 = 180
$ws.Cells.Item(169, 11).Value2 = 12000
$ws.Cells.Item(169, 12).Value2 = 14000
$ws.Cells.Item(169, 13).Value2 = 13111
$ws.Cells.Item(169, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(169, 16).Value2 = 219
